$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 8718.205
$ws.Range("I62").Value = 9373.030000000001
$ws.Range("J62").Value = 5116.6665
$ws.Range("K62").Value = 9373.030000000001
$ws.Range("L62").Value = 5116.6665
$ws.Range("M62").Value = -8749.030000000001
$ws.Range("N62").Value = -6364.6665

$ws.Range("H65").Value = 8718.205
$ws.Range("I65").Value = 9373.030000000001
$ws.Range("J65").Value = 5116.6665
$ws.Range("K65").Value = 46865.15
$ws.Range("L65").Value = 25583.3325
$ws.Range("M65").Value = -43745.15
$ws.Range("N65").Value = -31823.3325

$ws.Range("H137").Value = 46797
$ws.Range("J137").Value = 1390.7273
$ws.Range("L137").Value = 4172.1819
$ws.Range("N137").Value = -9272.1819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 7682.5
$ws.Range("I2").Value = 670
$ws.Range("J2").Value = 25213.75
$ws.Range("K2").Value = 670
$ws.Range("L2").Value = 25213.75
$ws.Range("M2").Value = -557
$ws.Range("N2").Value = -25439.75

$ws.Range("H61").Value = 1647.6
$ws.Range("I61").Value = 1397.1765
$ws.Range("K61").Value = 1397.1765
$ws.Range("M61").Value = -1185.1765

$ws.Range("H116").Value = 7682.5
$ws.Range("I116").Value = 670
$ws.Range("J116").Value = 25213.75
$ws.Range("K116").Value = 670
$ws.Range("L116").Value = 25213.75
$ws.Range("M116").Value = 1624
$ws.Range("N116").Value = -29801.75

$ws.Range("H132").Value = 1809473.8
$ws.Range("I132").Value = 2128241
$ws.Range("K132").Value = 6384723
$ws.Range("M132").Value = -6382193

$ws.Range("H136").Value = 1647.6
$ws.Range("I136").Value = 1397.1765
$ws.Range("K136").Value = 4191.529500000001
$ws.Range("M136").Value = -1641.529500000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 7682.5
$ws.Range("I3").Value = 670
$ws.Range("J3").Value = 25213.75
$ws.Range("K3").Value = 670
$ws.Range("L3").Value = 25213.75
$ws.Range("M3").Value = -556
$ws.Range("N3").Value = -25441.75

$ws.Range("H94").Value = 1448
$ws.Range("I94").Value = 472.16666
$ws.Range("J94").Value = 3399.6667
$ws.Range("K94").Value = 472.16666
$ws.Range("L94").Value = 3399.6667
$ws.Range("M94").Value = -21.16665999999998
$ws.Range("N94").Value = -4301.6667

$ws.Range("H134").Value = 55139.094
$ws.Range("I134").Value = 2577.0715
$ws.Range("K134").Value = 7731.2145
$ws.Range("M134").Value = -5196.2145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1800.9231
$ws.Range("I58").Value = 1374.909
$ws.Range("J58").Value = 2113.3333
$ws.Range("K58").Value = 1374.909
$ws.Range("L58").Value = 2113.3333
$ws.Range("M58").Value = -1171.909
$ws.Range("N58").Value = -2519.3333

$ws.Range("H107").Value = 473.0645
$ws.Range("I107").Value = 366.75
$ws.Range("J107").Value = 586.4666999999999
$ws.Range("K107").Value = 366.75
$ws.Range("L107").Value = 586.4666999999999
$ws.Range("M107").Value = 1553.25
$ws.Range("N107").Value = -4426.4667

$ws.Range("H136").Value = 1800.9231
$ws.Range("I136").Value = 1374.909
$ws.Range("J136").Value = 2113.3333
$ws.Range("K136").Value = 4124.727000000001
$ws.Range("L136").Value = 6339.999899999999
$ws.Range("M136").Value = -1574.727000000001
$ws.Range("N136").Value = -11439.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 5048.154
$ws.Range("I75").Value = 3266.5
$ws.Range("J75").Value = 5372.091
$ws.Range("K75").Value = 9799.5
$ws.Range("L75").Value = 16116.273
$ws.Range("M75").Value = -8801.5
$ws.Range("N75").Value = -18112.273

$ws.Range("H78").Value = 5048.154
$ws.Range("I78").Value = 3266.5
$ws.Range("J78").Value = 5372.091
$ws.Range("K78").Value = 29398.5
$ws.Range("L78").Value = 48348.819
$ws.Range("M78").Value = -24406.5
$ws.Range("N78").Value = -58332.819

$ws.Range("H87").Value = 46461.07
$ws.Range("I87").Value = 3080
$ws.Range("J87").Value = 55891.74
$ws.Range("K87").Value = 9240
$ws.Range("L87").Value = 167675.22
$ws.Range("M87").Value = -7992
$ws.Range("N87").Value = -170171.22

$ws.Range("H90").Value = 46461.07
$ws.Range("I90").Value = 3080
$ws.Range("J90").Value = 55891.74
$ws.Range("K90").Value = 27720
$ws.Range("L90").Value = 503025.66
$ws.Range("M90").Value = -21480
$ws.Range("N90").Value = -515505.66

$ws.Range("H131").Value = 46875720
$ws.Range("J131").Value = 50000730
$ws.Range("L131").Value = 150002190
$ws.Range("N131").Value = -150012270

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 15955.777
$ws.Range("I102").Value = 6999.5884
$ws.Range("J102").Value = 31181.3
$ws.Range("K102").Value = 6999.5884
$ws.Range("L102").Value = 31181.3
$ws.Range("M102").Value = -5377.5884
$ws.Range("N102").Value = -34425.3

$ws.Range("H132").Value = 86278.414
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 86278.414
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 258835.242
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -263895.242

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1602.5869
$ws.Range("I93").Value = 1379.6333
$ws.Range("J93").Value = 2020.625
$ws.Range("K93").Value = 1379.6333
$ws.Range("L93").Value = 2020.625
$ws.Range("M93").Value = -131.6333
$ws.Range("N93").Value = -4516.625

$ws.Range("H132").Value = 434300.1
$ws.Range("I132").Value = 205671.1
$ws.Range("K132").Value = 617013.3
$ws.Range("M132").Value = -614483.3

$ws.Range("H136").Value = 359336.72
$ws.Range("I136").Value = 716103.4399999999
$ws.Range("J136").Value = 2570
$ws.Range("K136").Value = 2148310.32
$ws.Range("L136").Value = 7710
$ws.Range("M136").Value = -2145760.32
$ws.Range("N136").Value = -12810

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 229.25
$ws.Range("I113").Value = 217.78947
$ws.Range("J113").Value = 253.44444
$ws.Range("K113").Value = 653.36841
$ws.Range("L113").Value = 760.33332
$ws.Range("M113").Value = 1516.63159
$ws.Range("N113").Value = -5100.33332

$ws.Range("H122").Value = 2909.4
$ws.Range("I122").Value = 1502.091
$ws.Range("J122").Value = 4629.4443
$ws.Range("K122").Value = 4506.272999999999
$ws.Range("L122").Value = 13888.3329
$ws.Range("M122").Value = -2056.272999999999
$ws.Range("N122").Value = -18788.3329

$ws.Range("H126").Value = 614.7353000000001
$ws.Range("I126").Value = 579.5172
$ws.Range("J126").Value = 819
$ws.Range("K126").Value = 1738.5516
$ws.Range("L126").Value = 2457
$ws.Range("M126").Value = 731.4484
$ws.Range("N126").Value = -7397

$ws.Range("H132").Value = 9203.429
$ws.Range("I132").Value = 1257.8334
$ws.Range("J132").Value = 15162.625
$ws.Range("K132").Value = 3773.5002
$ws.Range("L132").Value = 45487.875
$ws.Range("M132").Value = -1243.5002
$ws.Range("N132").Value = -50547.875

$ws.Range("H136").Value = 3394317.5
$ws.Range("I136").Value = 4763080.5
$ws.Range("J136").Value = 1113045.5
$ws.Range("K136").Value = 14289241.5
$ws.Range("L136").Value = 3339136.5
$ws.Range("M136").Value = -14286691.5
$ws.Range("N136").Value = -3344236.5
